$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change) per source diff.
# D-column price cells are plain text in the source file (e.g. "26.004.56" or "0.9989").
# Force a Text number format before assigning so Excel does not reinterpret them as numbers,
# then restore the default "Normal" style so the cell formatting matches the original file.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.996.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.751.90"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5207"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2852"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.749.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07029"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6473"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.536"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9983"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9987"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.997.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006638"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.976.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.159"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.671"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.173"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.504"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.855"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08327"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.664"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.447"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04453"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9899"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6113"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01591"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.954"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9984"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3887"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7362"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.033"
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.377"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1119"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.608"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
